$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "FLAGYL 500MG 20 TAB." product row (row 38) entirely; rows below shift up.
$ws.Rows.Item(38).Delete()

# Update the grand-total cell (now shifted from P91 to P90) to reflect the removed row's price.
$ws.Range("P90").Value = 3864.9949999999999

# Update the generated timestamp footer (now shifted from A92 to A91).
$ws.Range("A91").Value = "Thursday, 14 August, 2025 5:34 PM"
